$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Rows 45-47: the "bsecode" column (E) was stored as text; convert to a
# real numeric value (matches F/G/H columns' numeric typing).
$ws.Cells.Item(45, 5).Value = 20
$ws.Cells.Item(46, 5).Value = 531344
$ws.Cells.Item(47, 5).Value = 505537

# New row 48: another screener hit appended to the bottom of the table.
$ws.Cells.Item(48, 1).Value = "25/06/2024 04:44:42"
$ws.Cells.Item(48, 2).Value = 1
$ws.Cells.Item(48, 3).Value = "FACT"
$ws.Cells.Item(48, 4).Value = "Fertilizers And Chemicals Travancore Limited"

# bsecode on the freshly-added row keeps the (pre-fix) text representation,
# just like E45:E47 originally had before this edit cleaned them up.
$ws.Cells.Item(48, 5).Value = "'590024"
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(48, 6).Value = -1.11
$ws.Cells.Item(48, 7).Value = 1009.7
$ws.Cells.Item(48, 8).Value = 1109622
